# The workbook gained a brand-new weekly data point. In the OOXML this shows
# up as a newly inserted row 55 (pushing the previously-existing rows 55..158
# down to 56..159, all other columns staying identical), so we replicate that
# with a real row insert followed by populating the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 55, shifting rows 55-158 down
# to 56-159 (and extending the sheet's used range to row 159).
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new observation's data.
$ws.Cells.Item(55, 1).Value  = 8
$ws.Cells.Item(55, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(55, 3).Value  = "Coquimbo"
$ws.Cells.Item(55, 4).Value  = 44775
$ws.Cells.Item(55, 5).Value  = 4
$ws.Cells.Item(55, 6).Value  = 100112044
$ws.Cells.Item(55, 7).Value  = "Perejil"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 2400
$ws.Cells.Item(55, 11).Value = 2000
$ws.Cells.Item(55, 12).Value = 2500
$ws.Cells.Item(55, 13).Value = 2250
$ws.Cells.Item(55, 14).Value = '$/atado 1 a 1,5 kilos'
$ws.Cells.Item(55, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(55, 16).Value = 1500
$ws.Cells.Item(55, 17).Value = 1.5
$ws.Cells.Item(55, 18).Value = "Hortaliza"
